$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1510.6666
$ws.Range("I40").Value = 1524.75
$ws.Range("J40").Value = 1499.4
$ws.Range("K40").Value = 1524.75
$ws.Range("L40").Value = 1499.4
$ws.Range("M40").Value = -1349.75
$ws.Range("N40").Value = -1849.4
$ws.Range("H64").Value = 2904.7144
$ws.Range("J64").Value = 2999.9092
$ws.Range("L64").Value = 2999.9092
$ws.Range("N64").Value = -3495.9092
$ws.Range("H67").Value = 2904.7144
$ws.Range("J67").Value = 2999.9092
$ws.Range("L67").Value = 2999.9092
$ws.Range("N67").Value = -4715.9092
$ws.Range("H76").Value = 5053463
$ws.Range("H79").Value = 5053463
$ws.Range("H132").Value = 4050.5454
$ws.Range("I132").Value = 3890.65
$ws.Range("K132").Value = 11671.95
$ws.Range("M132").Value = -9141.950000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3229.36
$ws.Range("I61").Value = 1812.56
$ws.Range("K61").Value = 1812.56
$ws.Range("M61").Value = -1600.56
$ws.Range("H74").Value = 11366121
$ws.Range("I74").Value = 1503.5834
$ws.Range("K74").Value = 1503.5834
$ws.Range("M74").Value = -629.5834
$ws.Range("H77").Value = 11366121
$ws.Range("I77").Value = 1503.5834
$ws.Range("K77").Value = 7517.916999999999
$ws.Range("M77").Value = -3149.916999999999
$ws.Range("H88").Value = 2289.5557
$ws.Range("I88").Value = 1826.5
$ws.Range("J88").Value = 2660
$ws.Range("K88").Value = 1826.5
$ws.Range("L88").Value = 2660
$ws.Range("M88").Value = -1420.5
$ws.Range("N88").Value = -3472
$ws.Range("H91").Value = 2289.5557
$ws.Range("I91").Value = 1826.5
$ws.Range("J91").Value = 2660
$ws.Range("K91").Value = 1826.5
$ws.Range("L91").Value = 2660
$ws.Range("M91").Value = -422.5
$ws.Range("N91").Value = -5468
$ws.Range("H136").Value = 3229.36
$ws.Range("I136").Value = 1812.56
$ws.Range("K136").Value = 5437.68
$ws.Range("M136").Value = -2887.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 3000
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 859.625
$ws.Range("I10").Value = 859.625
$ws.Range("K10").Value = 859.625
$ws.Range("M10").Value = -720.625
$ws.Range("H14").Value = 850
$ws.Range("I14").Value = 850
$ws.Range("K14").Value = 850
$ws.Range("M14").Value = -680
$ws.Range("H62").Value = 4482.456
$ws.Range("I62").Value = 4482.456
$ws.Range("K62").Value = 4482.456
$ws.Range("M62").Value = -3858.456
$ws.Range("H65").Value = 4482.456
$ws.Range("I65").Value = 4482.456
$ws.Range("K65").Value = 22412.28
$ws.Range("M65").Value = -19292.28
$ws.Range("H132").Value = 2161
$ws.Range("I132").Value = 1895.4706
$ws.Range("K132").Value = 5686.4118
$ws.Range("M132").Value = -3156.4118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1216.2142
$ws.Range("I68").Value = 1200.4
$ws.Range("J68").Value = 1225
$ws.Range("K68").Value = 3601.2
$ws.Range("L68").Value = 3675
$ws.Range("M68").Value = -2790.2
$ws.Range("N68").Value = -5297
$ws.Range("H71").Value = 1216.2142
$ws.Range("I71").Value = 1200.4
$ws.Range("J71").Value = 1225
$ws.Range("K71").Value = 10803.6
$ws.Range("L71").Value = 11025
$ws.Range("M71").Value = -6747.6
$ws.Range("N71").Value = -19137
$ws.Range("H137").Value = 26977.488
$ws.Range("I137").Value = 5966.391
$ws.Range("J137").Value = 48943.637
$ws.Range("K137").Value = 17899.173
$ws.Range("L137").Value = 146830.911
$ws.Range("M137").Value = -12799.173
$ws.Range("N137").Value = -157030.911
$ws.Range("H140").Value = 1701.5333
$ws.Range("I140").Value = 1499.2858
$ws.Range("J140").Value = 4533
$ws.Range("K140").Value = 4497.857400000001
$ws.Range("L140").Value = 13599
$ws.Range("M140").Value = 682.1425999999992
$ws.Range("N140").Value = -23959

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13564.182
$ws.Range("I70").Value = 15289.556
$ws.Range("K70").Value = 15289.556
$ws.Range("M70").Value = -15019.556
$ws.Range("H73").Value = 13564.182
$ws.Range("I73").Value = 15289.556
$ws.Range("K73").Value = 15289.556
$ws.Range("M73").Value = -14353.556
$ws.Range("H80").Value = 820596.8
$ws.Range("I80").Value = 1502092.1
$ws.Range("J80").Value = 2802.4
$ws.Range("K80").Value = 1502092.1
$ws.Range("L80").Value = 2802.4
$ws.Range("M80").Value = -1501094.1
$ws.Range("N80").Value = -4798.4
$ws.Range("H83").Value = 820596.8
$ws.Range("I83").Value = 1502092.1
$ws.Range("J83").Value = 2802.4
$ws.Range("K83").Value = 7510460.5
$ws.Range("L83").Value = 14012
$ws.Range("M83").Value = -7505468.5
$ws.Range("N83").Value = -23996
$ws.Range("H122").Value = 13277.777
$ws.Range("I122").Value = 27375
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 82125
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -79675
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 2645.3333
$ws.Range("I132").Value = 1858.3529
$ws.Range("J132").Value = 3674.4614
$ws.Range("K132").Value = 5575.0587
$ws.Range("L132").Value = 11023.3842
$ws.Range("M132").Value = -3045.0587
$ws.Range("N132").Value = -16083.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2975.8333
$ws.Range("I40").Value = 2771
$ws.Range("K40").Value = 2771
$ws.Range("M40").Value = -2635
$ws.Range("H82").Value = 2184.889
$ws.Range("I82").Value = 1988.25
$ws.Range("J82").Value = 2342.2
$ws.Range("K82").Value = 1988.25
$ws.Range("L82").Value = 2342.2
$ws.Range("M82").Value = -1627.25
$ws.Range("N82").Value = -3064.2
$ws.Range("H85").Value = 2184.889
$ws.Range("I85").Value = 1988.25
$ws.Range("J85").Value = 2342.2
$ws.Range("K85").Value = 1988.25
$ws.Range("L85").Value = 2342.2
$ws.Range("M85").Value = -740.25
$ws.Range("N85").Value = -4838.2
$ws.Range("H122").Value = 4691.25
$ws.Range("I122").Value = 4844.4443
$ws.Range("J122").Value = 4494.2856
$ws.Range("K122").Value = 14533.3329
$ws.Range("L122").Value = 13482.8568
$ws.Range("M122").Value = -12083.3329
$ws.Range("N122").Value = -18382.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 5119221.5
$ws.Range("I132").Value = 2063.9143
$ws.Range("K132").Value = 6191.742899999999
$ws.Range("M132").Value = -3661.742899999999
